# Applies Betfair Back/Lay odds updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U3").Value = 2.4
$ws.Range("G4").Value = 1.64
$ws.Range("I4").Value = 7
$ws.Range("K4").Value = 4.6
$ws.Range("P4").Value = 2.06
$ws.Range("T4").Value = 1.94
$ws.Range("V4").Value = 1.17
$ws.Range("W4").Value = 2.54
$ws.Range("Q6").Value = 1.94
$ws.Range("O7").Value = 1.4
$ws.Range("I8").Value = 8.4
$ws.Range("P8").Value = 2.88
$ws.Range("T8").Value = 1.72
$ws.Range("U8").Value = 2.2
$ws.Range("W8").Value = 3.3
$ws.Range("J9").Value = 3.1
$ws.Range("O9").Value = 1.26
$ws.Range("S9").Value = 3.05
$ws.Range("T9").Value = 1.58
$ws.Range("U9").Value = 2
$ws.Range("O10").Value = 1.5
$ws.Range("H11").Value = 1.79
$ws.Range("G12").Value = 2.26
$ws.Range("P12").Value = 1.72
$ws.Range("S12").Value = 3.95
$ws.Range("T12").Value = 1.8
$ws.Range("V12").Value = 1.27
$ws.Range("W12").Value = 1.81
$ws.Range("F13").Value = 2.34
$ws.Range("G13").Value = 2.66
$ws.Range("I13").Value = 3.95
$ws.Range("O13").Value = 1.42
$ws.Range("Q13").Value = 2.24
$ws.Range("S13").Value = 3.9
$ws.Range("T13").Value = 1.9
$ws.Range("U13").Value = 1.92
$ws.Range("V13").Value = 1.34
$ws.Range("W13").Value = 1.6
$ws.Range("X13").Value = 12
$ws.Range("AG13").Value = 14
$ws.Range("H14").Value = 1.82
$ws.Range("Q14").Value = 1.7
$ws.Range("S14").Value = 1.7
$ws.Range("G15").Value = 2.32
$ws.Range("I15").Value = 4.3
$ws.Range("O15").Value = 1.45
$ws.Range("T15").Value = 1.98
$ws.Range("W15").Value = 1.76
$ws.Range("F16").Value = 1.92
$ws.Range("S16").Value = 1.84
$ws.Range("W16").Value = 1.93
$ws.Range("Z16").Value = 36
$ws.Range("F17").Value = 3.75
$ws.Range("H17").Value = 1.86
$ws.Range("I17").Value = 2.02
$ws.Range("L17").Value = 1.28
$ws.Range("M17").Value = 1.01
$ws.Range("P17").Value = 2.46
$ws.Range("R17").Value = 1.58
$ws.Range("U17").Value = 2.46
$ws.Range("V17").Value = 1.98
$ws.Range("W17").Value = 1.27
$ws.Range("F18").Value = 1.81
$ws.Range("G18").Value = 2.02
$ws.Range("H18").Value = 4.7
$ws.Range("O18").Value = 1.4
$ws.Range("W18").Value = 1.99
$ws.Range("Z18").Value = 46
$ws.Range("AE18").Value = 95
$ws.Range("AI18").Value = 110
$ws.Range("AJ18").Value = 25
$ws.Range("AM18").Value = 180
$ws.Range("F19").Value = 2.48
$ws.Range("G19").Value = 2.6
$ws.Range("H19").Value = 3.35
$ws.Range("L19").Value = 1.57
$ws.Range("U19").Value = 1.83
$ws.Range("W19").Value = 1.62
$ws.Range("F20").Value = 2.68
$ws.Range("G20").Value = 2.82
$ws.Range("H20").Value = 2.8
$ws.Range("J20").Value = 3.45
$ws.Range("P20").Value = 1.72
$ws.Range("Q20").Value = 2.18
$ws.Range("S20").Value = 3.9
$ws.Range("T20").Value = 1.83
$ws.Range("U20").Value = 1.97
$ws.Range("T21").Value = 2
$ws.Range("Z21").Value = 18
$ws.Range("AB21").Value = 9
$ws.Range("AL21").Value = 65
$ws.Range("Q22").Value = 1.53
$ws.Range("S22").Value = 2.34
$ws.Range("U22").Value = 2.76
$ws.Range("AC22").Value = 10.5
$ws.Range("P23").Value = 1.81
$ws.Range("X23").Value = 11.5
$ws.Range("G24").Value = 3.05
$ws.Range("H24").Value = 2.6
$ws.Range("I25").Value = 3.95
$ws.Range("P25").Value = 1.66
